$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The daily price sheet gets a brand-new "today" row inserted above the
# previous top row. Everything else shifts down by one row, which naturally
# pushes the previous last row (249 -> 250) without losing any data (the
# source table appears to pad the bottom with a repeat of the oldest known
# entry when it runs out of older history).
# ---------------------------------------------------------------------------

# 1) Make room: insert a blank row at row 2, shifting rows 2..249 down to
#    3..250. This also extends the used range to A1:F250 automatically.
$ws.Rows.Item(2).Insert()

# 2) The brand-new row reuses the same Basic Price / Circular Date / Circular
#    Link as the (now shifted-down) previous top row at row 3, since the
#    price had not changed yet on the new date - only the Date column is new.
$ws.Range("B2").Value = $ws.Range("B3").Value()
$ws.Range("C2").Value = $ws.Range("C3").Value()
$ws.Range("D2").Value = $ws.Range("D3").Value()
$ws.Range("E2").Value = $ws.Range("E3").Value()
$ws.Range("F2").Value = $ws.Range("F3").Value()
$ws.Range("A2").Value = "15-02-2026"

# 3) Row 2 was created with the header row's formatting (Insert copies the
#    format of the row above). Restore the normal body formatting - centered
#    alignment for every column, plus the numeric "0.000" format on D - to
#    match the rest of the table (row 3 is a perfect template).
$ws.Range("A3:F3").Copy()
$ws.Range("A2:F2").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# 4) Hyperlinks in this engine stay pinned to their original absolute cell
#    address when rows are inserted (they do not ride along with the shifted
#    cell content), so every Circular Link hyperlink needs to be rebuilt from
#    scratch against the final (now-correct) text in column F.
$ws.Hyperlinks.Delete()
$lastRow = $ws.Range("A1").End(-4121).Row()
for ($r = 2; $r -le $lastRow; $r++) {
    $target = $ws.Cells.Item($r, 6).Value()
    $ws.Hyperlinks.Add($ws.Cells.Item($r, 6), $target)
}
